# Modif excel pour lier heures a feuille de route, pour issu #7
#
# Link the "Ceremonies protocolaires" / "Souper" start times on the MAR
# sheet to the shuttle's expected-arrival time (E17, itself pulled from
# Navettes!E2) instead of leaving them as static/unlinked cells, and add
# a reminder note on the Navettes sheet that its data must be refreshed
# by hand when the itinerary changes.

$wb = $excel.ActiveWorkbook

# --- MAR sheet: tie schedule times to the shuttle arrival time (E17) ---
$wsMar = $wb.Worksheets.Item("MAR")

# Ceremonies protocolaires: 20 minutes after the shuttle's arrival
$wsMar.Range("F18").Value = 0.013888888888888888
$wsMar.Range("E18").Formula = "=E17+F18"

# Souper: 10 minutes after the shuttle's arrival
$wsMar.Range("F19").Value = 0.0069444444444444441
$wsMar.Range("E19").Formula = "=E17+F19"

# --- Navettes sheet: add a manual-refresh reminder note below the existing one ---
$wsNav = $wb.Worksheets.Item("Navettes")
$wsNav.Range("A12").Copy()
$wsNav.Range("A14").PasteSpecial(-4122)  # xlPasteFormats - reuse A12's note styling
$wsNav.Range("A14").Value = "À MAJ en actualisant les données manuellement lors de changement de l'itinéraire."

# --- Update the active sheet/selection state ---
$wsNav.Activate()
$wsNav.Range("A15").Select()

$wsMar.Activate()
$wsMar.Range("E23").Select()
